$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.397.47"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "3.902.50"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'602.13"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'171.56"
$ws.Range("D7").Value = "3.904.07"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("D11").Value = "'6.43"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "'0.0000262"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").Value = "'37.33"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "4.558.60"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "3.903.31"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "68.467.21"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "'18.20"
$ws.Range("E18").Value = "  +4.79%  "
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").Value = "'0.111"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'10.83"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "'470.25"
$ws.Range("E22").Value = "  -4.42%  "
$ws.Range("D23").Value = "'0.741"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").Value = "'0.0000164"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'9.95"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").Value = "'2.97"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").Value = "4.054.05"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "'7.80"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "'31.42"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").Value = "'9.43"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "3.873.16"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").Value = "'3.69"
$ws.Range("E38").Value = "  +14.53%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +12.29%  "
$ws.Range("D44").Value = "'0.313"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.99"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'424.96"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "'47.19"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").Value = "'27.48"
$ws.Range("E50").Value = "  +5.46%  "
$ws.Range("D51").Value = "'143.26"
$ws.Range("E51").Value = "  +0.27%  "
